$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview!G2 and de-de!H2 share the same original value "2016-09-05 17:13:41"
# and both move to "2016-09-05 17:14:28"
$wsOverview.Range("G2").Value = "2016-09-05 17:14:28"
$wsDeDe.Range("H2").Value = "2016-09-05 17:14:28"

# zh-cn!H2: Correspond Handoff Datetime 17:13:36 -> 17:14:23
$wsZhCn.Range("H2").Value = "2016-09-05 17:14:23"

# zh-cn!K2: Correspond Handback DateTime 17:13:54 -> 17:15:11
$wsZhCn.Range("K2").Value = "2016-09-05 17:15:11"

# de-de!K2: Correspond Handback DateTime 17:14:01 -> 17:15:22
$wsDeDe.Range("K2").Value = "2016-09-05 17:15:22"
